$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 5.611999999999999
$ws.Range("C3").Value = -12.11
$ws.Range("E3").Value = 16.548
$ws.Range("E6").Value = 16.864
$ws.Range("D8").Value = -8.49
$ws.Range("D11").Value = -7.394
$ws.Range("A12").Value = -21.629
$ws.Range("B14").Value = 6.114
$ws.Range("D14").Value = -7.678
$ws.Range("D15").Value = -8.175000000000001
$ws.Range("B26").Value = 6.488
$ws.Range("E27").Value = 16.691
$ws.Range("C30").Value = -12.413
$ws.Range("B31").Value = 6.371
$ws.Range("A32").Value = -21.385
$ws.Range("E33").Value = 17.558
$ws.Range("B35").Value = 8.020999999999999
$ws.Range("A36").Value = -21.044
$ws.Range("D36").Value = -8.169
$ws.Range("B37").Value = 8.260000000000002
$ws.Range("A38").Value = -20.223
$ws.Range("E39").Value = 16.092
$ws.Range("C44").Value = -12.289
$ws.Range("B45").Value = 5.389
$ws.Range("A46").Value = -21.515
$ws.Range("E47").Value = 16.615
$ws.Range("A54").Value = -21.894
$ws.Range("E54").Value = 16.77
$ws.Range("A55").Value = -22.082
$ws.Range("E56").Value = 16.716
$ws.Range("B57").Value = 6.090000000000001
$ws.Range("C58").Value = -12.813
$ws.Range("E58").Value = 16.691
$ws.Range("D64").Value = -7.81
$ws.Range("E66").Value = 17.482
$ws.Range("A67").Value = -21.588
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.497
$ws.Range("E72").Value = 17.033
$ws.Range("E82").Value = 17.002
$ws.Range("E83").Value = 16.851
$ws.Range("C84").Value = -12.781
$ws.Range("C89").Value = -12.339
$ws.Range("D89").Value = -7.487
$ws.Range("A91").Value = -21.632
$ws.Range("C91").Value = -11.069
$ws.Range("C92").Value = -11.864
$ws.Range("A99").Value = -20.682
$ws.Range("B100").Value = 6.006
$ws.Range("B102").Value = 7.529000000000001
$ws.Range("C102").Value = -12.327
